$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row (50) down to the new row (51)
# so number formats / styles match exactly (date format on A, 0.0000 on F).
[void]$ws.Range("A50:H50").Copy()
[void]$ws.Range("A51:H51").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data row (2025-02-21) appended to the return series.
$ws.Range("A51").Value = 45709
$ws.Range("B51").Value = -0.0171
$ws.Range("C51").Value = -0.016
$ws.Range("D51").Formula = "=B51+1"
$ws.Range("E51").Formula = "=1+C51"
$ws.Range("F51").Formula = "=C51-B51"
$ws.Range("G51").Formula = "=G50*D50"
$ws.Range("H51").Formula = "=H50*E50"

# Reflect the new extent of the data in the current selection, like Excel
# does when the used range grows and the whole table is (re)selected.
[void]$ws.Range("A1:H51").Select()
